# Daily attendance processing - reorder "Recorded By" list so the last
# recorder in the list is moved to the front (most-recent-first ordering),
# leaving single-value cells and the specific "System, admin@admin.com"
# pairing untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $val = $cell.Text

    if ([string]::IsNullOrEmpty($val)) {
        continue
    }

    if ($val -eq "System, admin@admin.com") {
        continue
    }

    $parts = $val -split ", "
    if ($parts.Count -lt 2) {
        continue
    }

    $n = $parts.Count
    $last = $parts[$n - 1]
    $rest = $parts[0..($n - 2)]
    $newParts = @($last) + $rest
    $newVal = $newParts -join ", "

    $cell.Value = $newVal
}
